$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 304, shifting existing rows 304:330 down to 305:331.
$ws.Rows.Item(304).Insert()

# Populate the newly inserted row 304 with the new record.
$ws.Cells.Item(304, 1).Value = 10
$ws.Cells.Item(304, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(304, 3).Value = "La Araucanía"
$ws.Cells.Item(304, 4).Value = 45166
$ws.Cells.Item(304, 5).Value = 9
$ws.Cells.Item(304, 6).Value = 100112005
$ws.Cells.Item(304, 7).Value = "Puerro"
$ws.Cells.Item(304, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(304, 9).Value = "Primera"
$ws.Cells.Item(304, 10).Value = 80
$ws.Cells.Item(304, 11).Value = 9000
$ws.Cells.Item(304, 12).Value = 9000
$ws.Cells.Item(304, 13).Value = 9000
$ws.Cells.Item(304, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(304, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(304, 16).Value = 750
$ws.Cells.Item(304, 17).Value = 12
$ws.Cells.Item(304, 18).Value = "Hortaliza"
